# Weekly price-list update: a new weekly observation row is inserted at
# row 7 (pushing the existing rows 7-42 down to 8-43), and populated with
# this week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 7 - shifts rows 7:42 down to 8:43
$ws.Rows("7:7").Insert()

# Fill the newly inserted row 7 with the new weekly record
$ws.Cells.Item(7, 1).Value  = 9
$ws.Cells.Item(7, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(7, 3).Value  = "Metropolitana"
$ws.Cells.Item(7, 4).Value  = 44545
$ws.Cells.Item(7, 5).Value  = 13
$ws.Cells.Item(7, 6).Value  = 100112029
$ws.Cells.Item(7, 7).Value  = "Orégano"
$ws.Cells.Item(7, 8).Value  = "Sin especificar"
$ws.Cells.Item(7, 9).Value  = "Primera"
$ws.Cells.Item(7, 10).Value = 25
$ws.Cells.Item(7, 11).Value = 9000
$ws.Cells.Item(7, 12).Value = 10000
$ws.Cells.Item(7, 13).Value = 9480
$ws.Cells.Item(7, 14).Value = "`$/docena de atados"
$ws.Cells.Item(7, 15).Value = "Región Metropolitana"
$ws.Cells.Item(7, 16).Value = 3160
$ws.Cells.Item(7, 17).Value = 3
$ws.Cells.Item(7, 18).Value = "Hortaliza"
